$d = $word.ActiveDocument

# Find the paragraph that contains only "Edison Achalma" styled as Author,
# which appears right after the main title (Heading1). We scan the
# Paragraphs collection to locate it precisely, then insert a new
# paragraph right after it with the same "Author" style.

$found = $false
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    $styleName = $p.Range.Style.NameLocal
    if (($text -eq "Edison Achalma`r") -and ($styleName -eq "Author")) {
        # Build a collapsed point range positioned exactly at the end of this
        # paragraph (i.e. right after its paragraph mark / at the start of
        # the following paragraph) and insert a new paragraph mark there.
        $endPos = $p.Range.End
        $insertRange = $d.Range($endPos, $endPos)
        $insertRange.InsertParagraphAfter()

        # The newly created paragraph now sits right after the original
        # "Edison Achalma" paragraph. Set its style and text.
        $newPara = $p.Next(1)
        $newPara.Range.Style = "Author"
        $newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"

        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not find the 'Edison Achalma' Author paragraph"
}
